$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39; existing rows 39-56 shift down to 40-57.
$ws.Rows.Item(39).Insert()

# Populate the newly inserted row 39 with the new weekly price record.
$ws.Cells.Item(39, 1).Value = 11
$ws.Cells.Item(39, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(39, 3).Value = "Bíobío"
$ws.Cells.Item(39, 4).Value = 44488
$ws.Cells.Item(39, 5).Value = 8
$ws.Cells.Item(39, 6).Value = 100112024
$ws.Cells.Item(39, 7).Value = "Choclo"
$ws.Cells.Item(39, 8).Value = "Dulce o Americano"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 100
$ws.Cells.Item(39, 11).Value = 44000
$ws.Cells.Item(39, 12).Value = 45000
$ws.Cells.Item(39, 13).Value = 44500
$ws.Cells.Item(39, 14).Value = "$/malla 70 unidades"
$ws.Cells.Item(39, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(39, 16).Value = 636
$ws.Cells.Item(39, 17).Value = 70
$ws.Cells.Item(39, 18).Value = "Hortaliza"

# Make sure the date column keeps the same date style used by the rest of column D.
$ws.Cells.Item(39, 4).NumberFormat = $ws.Cells.Item(40, 4).NumberFormat
